$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.1
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.63
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8
$ws.Range("AN2").Value = 4
$ws.Range("AZ2").Value = 23
$ws.Range("BA2").Value = 67
$ws.Range("BC2").Value = 126
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("G4").Value = 3.7
$ws.Range("I4").Value = 2.1
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 41
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 8.5
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 7
$ws.Range("AI4").Value = 9.5
$ws.Range("AK4").Value = 19
$ws.Range("AL4").Value = 17
$ws.Range("AN4").Value = 5.5
$ws.Range("AO4").Value = 21
$ws.Range("AR4").Value = 101
$ws.Range("AX4").Value = 4
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 2.2
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.91
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 10
$ws.Range("AF5").Value = 51
$ws.Range("AG5").Value = 251
$ws.Range("AH5").Value = 12
$ws.Range("AI5").Value = 21
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 10
$ws.Range("AP5").Value = 21
$ws.Range("AS5").Value = 151
$ws.Range("AT5").Value = 2.75
$ws.Range("AY5").Value = 21
$ws.Range("BB5").Value = 101
